$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.082.74"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.204.39"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'537.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'145.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("D9").Value = "'7.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("D11").Value = "'0.434"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.16%  "
$ws.Range("D12").Value = "3.754.57"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").Value = "'25.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "60.108.98"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "3.215.18"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "'6.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D20").Value = "'8.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'375.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").Value = "'8.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.11%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "0.0₃0895"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "'22.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "'6.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "'5.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").Value = "'6.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.42%  "
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").Value = "'156.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "2.802.59"
$ws.Range("E37").Value = "  +6.19%  "
$ws.Range("D38").Value = "'25.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").Value = "'39.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("D43").Value = "'0.0292"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.08%  "
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("D46").Value = "3.245.89"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'0.982"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").Value = "'0.813"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.38%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("E51").Value = "  -0.03%  "

Write-Output "Applied 81 cell updates"
